$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 106, shifting existing rows 106:196 down
# to 107:197 (xlShiftDown = -4121).
$ws.Rows.Item(106).Insert(-4121)

# Populate the new row 106 with the new market-report record.
$ws.Cells.Item(106, 1).Value  = 4
$ws.Cells.Item(106, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(106, 3).Value  = "Los Lagos"
$ws.Cells.Item(106, 4).Value  = 44589
$ws.Cells.Item(106, 5).Value  = 10
$ws.Cells.Item(106, 6).Value  = "Fruta"
$ws.Cells.Item(106, 7).Value  = 100108
$ws.Cells.Item(106, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(106, 9).Value  = 100108005
$ws.Cells.Item(106, 10).Value = "Piña"
$ws.Cells.Item(106, 11).Value = "Caramelo"
$ws.Cells.Item(106, 12).Value = "Tercera"
$ws.Cells.Item(106, 13).Value = 200
$ws.Cells.Item(106, 14).Value = 17000
$ws.Cells.Item(106, 15).Value = 18000
$ws.Cells.Item(106, 16).Value = 17500
$ws.Cells.Item(106, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(106, 18).Value = "Ecuador"
$ws.Cells.Item(106, 19).Value = 1094
$ws.Cells.Item(106, 20).Value = 16
